# Trade #40 (internal id 68) closed at 2026-02-18 00:15:07 - unknown UNKNOWN +0.000%
# Also records two brand-new OPEN trades (#97 momentum, #98 HighProbConvergence)
# that were generated in the same tick, and rolls the aggregate stats on the
# Summary / Strategy Status sheets forward to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - roll aggregate counters forward
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 68      # Total Trades
$summary.Range("B9").Value = 52.94   # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - momentum strategy row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D11").Value = 5    # Trades
$status.Range("G11").Value = 20   # Win Rate %

# ---------------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a. Close out trade #68 (row 69)
$allTrades.Cells.Item(69, 7).Value = 0.03          # G: Exit Price
$allTrades.Cells.Item(69, 8).Value = "CLOSED"      # H: Status
$allTrades.Cells.Item(69, 11).Value = 99.68000000000001  # K: Capital After
$allTrades.Cells.Item(69, 12).Value = "early_exit" # L: Exit Reason
$allTrades.Cells.Item(69, 13).Value = 0.12         # M: Duration (min)

# 3b. New row 98 - trade #97 (momentum, DOWN, still OPEN)
$allTrades.Cells.Item(98, 1).Value = 97
$allTrades.Cells.Item(98, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(98, 3).Value = "00:15:00"
$allTrades.Cells.Item(98, 4).Value = "momentum"
$allTrades.Cells.Item(98, 5).Value = "DOWN"
$allTrades.Cells.Item(98, 6).Value = 0.03
$allTrades.Cells.Item(98, 8).Value = "OPEN"
$allTrades.Cells.Item(98, 9).Value = 0
$allTrades.Cells.Item(98, 10).Value = 0
$allTrades.Cells.Item(98, 11).Value = 99.6787371310913
$allTrades.Cells.Item(98, 13).Value = 0
$allTrades.Cells.Item(98, 14).Value = 0
$allTrades.Cells.Item(98, 15).Value = 0
$allTrades.Cells.Item(98, 16).Value = 0.9
$allTrades.Cells.Item(98, 17).Value = "Downward momentum: -40.404% over 10 samples"

# 3c. New row 99 - trade #98 (HighProbConvergence, UP, still OPEN)
$allTrades.Cells.Item(99, 1).Value = 98
$allTrades.Cells.Item(99, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(99, 3).Value = "00:15:01"
$allTrades.Cells.Item(99, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(99, 5).Value = "UP"
$allTrades.Cells.Item(99, 6).Value = 0.98
$allTrades.Cells.Item(99, 8).Value = "OPEN"
$allTrades.Cells.Item(99, 9).Value = 0
$allTrades.Cells.Item(99, 10).Value = 0
$allTrades.Cells.Item(99, 11).Value = 100.0565626577805
$allTrades.Cells.Item(99, 13).Value = 0
$allTrades.Cells.Item(99, 14).Value = 0
$allTrades.Cells.Item(99, 15).Value = 0
$allTrades.Cells.Item(99, 16).Value = 0.95
$allTrades.Cells.Item(99, 17).Value = "Mean reversion UP: price 33.71% below mean (z=-1.73)"

# ---------------------------------------------------------------------------
# 4. momentum sheet (strategy-specific log)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

# 4a. Close out trade #68 (row 6)
$momentum.Cells.Item(6, 7).Value = 0.03                  # G: Exit Price
$momentum.Cells.Item(6, 8).Value = "CLOSED"              # H: Status
$momentum.Cells.Item(6, 11).Value = 99.68000000000001    # K: Capital After
$momentum.Cells.Item(6, 16).Value = "early_exit"         # P: Exit Reason
$momentum.Cells.Item(6, 17).Value = 0.12                 # Q: Duration (min)

# 4b. New row 22 - trade #97 (momentum, DOWN, still OPEN)
$momentum.Cells.Item(22, 1).Value = 97
$momentum.Cells.Item(22, 2).Value = "'2026-02-18"
$momentum.Cells.Item(22, 3).Value = "00:15:00"
$momentum.Cells.Item(22, 4).Value = "momentum"
$momentum.Cells.Item(22, 5).Value = "DOWN"
$momentum.Cells.Item(22, 6).Value = 0.03
$momentum.Cells.Item(22, 8).Value = "OPEN"
$momentum.Cells.Item(22, 9).Value = 0
$momentum.Cells.Item(22, 10).Value = 0
$momentum.Cells.Item(22, 11).Value = 99.6787371310913
$momentum.Cells.Item(22, 12).Value = 0
$momentum.Cells.Item(22, 13).Value = 0
$momentum.Cells.Item(22, 14).Value = 0.9
$momentum.Cells.Item(22, 15).Value = "Downward momentum: -40.404% over 10 samples"
$momentum.Cells.Item(22, 17).Value = 0

# ---------------------------------------------------------------------------
# 5. HighProbConvergence sheet (strategy-specific log)
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")

# 5a. New row 11 - trade #98 (HighProbConvergence, UP, still OPEN)
$hpc.Cells.Item(11, 1).Value = 98
$hpc.Cells.Item(11, 2).Value = "'2026-02-18"
$hpc.Cells.Item(11, 3).Value = "00:15:01"
$hpc.Cells.Item(11, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(11, 5).Value = "UP"
$hpc.Cells.Item(11, 6).Value = 0.98
$hpc.Cells.Item(11, 8).Value = "OPEN"
$hpc.Cells.Item(11, 9).Value = 0
$hpc.Cells.Item(11, 10).Value = 0
$hpc.Cells.Item(11, 11).Value = 100.0565626577805
$hpc.Cells.Item(11, 12).Value = 0
$hpc.Cells.Item(11, 13).Value = 0
$hpc.Cells.Item(11, 14).Value = 0.95
$hpc.Cells.Item(11, 15).Value = "Mean reversion UP: price 33.71% below mean (z=-1.73)"
$hpc.Cells.Item(11, 17).Value = 0
